$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds "K" (formerly "Strike#") values - regenerated per commit message.
$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 3
$ws.Range("G4").Value = 2
$ws.Range("G5").Value = 3
$ws.Range("G6").Value = 0
$ws.Range("G7").Value = 1
$ws.Range("G8").Value = 1
$ws.Range("G9").Value = 1
